# Update the Metrics sheet's source values (B2:B13). Everything that
# depends on these (the "today" sheet's B11:B22 / E11:E22 / F11:F22
# formulas, plus TODAY()-1 in A1) recalculates automatically.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metrics")
$ws.Range("B2").Value = 228112.46000000005
$ws.Range("B3").Value = 200752.49
$ws.Range("B4").Value = 70494.600000000006
$ws.Range("B5").Value = 9315
$ws.Range("B6").Value = 5024358.2100000018
$ws.Range("B7").Value = 4242829.1700000009
$ws.Range("B8").Value = 1477454.43
$ws.Range("B9").Value = 195522
$ws.Range("B10").Value = 33489739.20000001
$ws.Range("B11").Value = 31518104.329999998
$ws.Range("B12").Value = 11759176.470000001
$ws.Range("B13").Value = 1293152

# Update each sheet's remembered selection. "today" must be selected
# last so it stays the active/visible tab, matching the source workbook
# (workbook.xml activeTab points at "today", whose sheetView carries
# tabSelected="1").
[void]$ws.Range("F15").Select()

$today = $wb.Worksheets.Item("today")
[void]$today.Range("D6").Select()
